$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the vector3 values in D5 and D6
$ws.Range("D5").Value = "231146|18953|3034"
$ws.Range("D6").Value = "212355|18550.12|2758.65"

# Remove rich-text formatting from B5 (battal_world_nameLV_20) by rewriting as plain text
$ws.Range("B5").Value = "battal_world_nameLV_20"

# Update window size
$excel.ActiveWindow.Width = 24045
$excel.ActiveWindow.Height = 12375

# Update selection
$ws.Range("K10").Select()
